$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts C:I -> D:J)
$ws.Columns.Item(3).Insert()

# Set the new header cell value
$ws.Range("C1").Value = "Industry"

# Match formatting of the other header cells (e.g. B1) for the new header
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Restore the text since PasteSpecial(formats) shouldn't touch it, but ensure value stays correct
$ws.Range("C1").Value = "Industry"
